# Adv-diff code for CAS
# Applies the edits to the "adv_diff" worksheet: fills in a couple of new
# measurements, inserts two blank spacer rows, and moves the active
# sheet/selection so "adv_diff" (not "diff") is the tab shown when the
# workbook is reopened.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("adv_diff")

# --- Row 19: fill in the missing measurement values -----------------------
$ws.Cells.Item(19, 3).Value = 0.00000010523
$ws.Cells.Item(19, 10).Value = 0.000000082885000000000004

# --- Row 20: label it with the two new shared strings ----------------------
# Set J20 ("random") before C20 ("atl/arc patch") so the new shared-string
# table entries land in the same order as the target workbook.
$ws.Cells.Item(20, 10).Value = "random"
$ws.Cells.Item(20, 10).ClearFormats()
$ws.Cells.Item(20, 3).Value = "atl/arc patch"

# --- Insert two blank rows before the old row 21 ---------------------------
# This pushes the former rows 21-25 down to 23-27, matching the new
# dimension (A1:L27).
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(21).Insert()

# The inserted rows copied column J's style/value down from the old row 21;
# the target has no J cell at all on these two spacer rows, so remove it.
$ws.Cells.Item(21, 10).Clear()
$ws.Cells.Item(22, 10).Clear()

# --- Make "adv_diff" the active sheet/tab with the new selection -----------
$ws.Activate()
$ws.Range("A20").Select()
